$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the date stored in O2 (datum column) from 2018-02-01 (43132) to 2018-01-14 (43114)
$ws.Range("O2").Value = 43114

# Give that cell a custom date number format (yyyy/mm/dd;@) instead of the built-in one
$ws.Range("O2").NumberFormat = "yyyy/mm/dd;@"

# Set the page to portrait orientation
$ws.PageSetup.Orientation = 1

# Move the active selection from N3 to N2
$ws.Range("N2").Select()
